$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13's "edad"/age-like column (D) needs to hold the text value "15"
# (previously "21"), stored as text (not auto-converted to a number) just
# like the rest of the column. Temporarily mark the cell as Text so the
# numeric-looking value is kept as a literal string, then restore the
# cell's style to Normal so no visible formatting change is left behind.
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "15"
$cell.Style = "Normal"
